$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Merge the two runs that make up "Mon Sep 24" / " 11:26:28 PDT 2017"
#    into a single run "Mon Sep 24 11:26:28 PDT 2017" (this text is
#    unique in the document, so a plain Find & Replace is safe and also
#    naturally coalesces the two runs into one, which is what the diff
#    shows).
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Mon Sep 24 11:26:28 PDT 2017", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Mon Sep 24 11:26:28 PDT 2017", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Find the last "Amount Received mode ... - CASH" paragraph (the end
#    of the last purchase record) and insert a brand-new purchase
#    record right after it, in front of the trailing blank paragraphs.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Start = 0
$rng.End = $d.Content.End
$lastStart = -1
$lastEnd = -1
while ($rng.Find.Execute("- CASH", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
  $lastStart = $rng.Start
  $lastEnd = $rng.End
  $rng.Collapse(0)
  $rng.End = $d.Content.End
}

# $lastEnd is the position right after the "H" of "- CASH" and right
# before that paragraph's end-of-paragraph mark.  Moving one more
# position forward lands right after the paragraph mark - i.e. exactly
# where the new paragraphs need to be inserted.
$insPos = $lastEnd + 1
$insRng = $d.Range($insPos, $insPos)

function Tabs($n) {
  $s = ""
  for ($i = 0; $i -lt $n; $i++) {
    $s = $s + [char]9
  }
  return $s
}

$CR = [char]13
$dashes = "---------------------------------------------------------------"

$block = "Tue Sep 25 11:15:39 PDT 2017" + $CR +
         "Person Name" + (Tabs 4) + "- TRM" + $CR +
         $dashes + $CR +
         "Item Name" + (Tabs 4) + "- CHOWCHOW EVE" + $CR +
         "Number of Pockets" + (Tabs 3) + "- 1" + $CR +
         "Number of KGs" + (Tabs 3) + "- 64" + $CR +
         "Rate" + (Tabs 5) + "- 15" + $CR +
         "Total Price" + (Tabs 4) + "- 960.0" + $CR +
         "Amount balance" + (Tabs 3) + "- 37871.0" + $CR +
         "" + $CR

$insRng.InsertBefore($block)

# ---------------------------------------------------------------------
# 3) Make the freshly inserted "Amount balance ... - 37871.0" line bold
#    (matching the bold "Amount balance" lines used elsewhere in the
#    document for the final balance of a purchase record).  Using the
#    Paragraphs collection (rather than Find + MoveEnd) ensures the
#    paragraph mark itself is included, so both the run and the
#    paragraph mark pick up the bold formatting, just like the diff.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
  $p = $d.Paragraphs.Item($i)
  if ($p.Range.Start -ge $insPos -and $p.Range.Text.IndexOf("- 37871.0") -ge 0) {
    $p.Range.Font.Bold = 1
  }
}
